# Minor fix related to tables structure
# - Adds a "name" column/field to both the "hotel" and "room" tables.
# - Re-draws the table borders (every column now gets a medium left/right
#   edge, thin edges between data rows, medium edges at the outer
#   top/bottom of each section).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

$xlThin = 2
$xlMedium = -4138
$xlLineStyleNone = -4142

$xlCenter = -4108

function Set-CellBorders {
    param(
        $cell,
        $left,
        $right,
        $top,
        $bottom
    )
    $r = $ws.Range($cell)
    if ($left) { $r.Borders.Item($xlEdgeLeft).Weight = $left } else { $r.Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone }
    if ($right) { $r.Borders.Item($xlEdgeRight).Weight = $right } else { $r.Borders.Item($xlEdgeRight).LineStyle = $xlLineStyleNone }
    if ($top) { $r.Borders.Item($xlEdgeTop).Weight = $top } else { $r.Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone }
    if ($bottom) { $r.Borders.Item($xlEdgeBottom).Weight = $bottom } else { $r.Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone }
}

function Format-Cell {
    param(
        $cell,
        $value,
        $bold = $false,
        $size = 11
    )
    $r = $ws.Range($cell)
    $r.Value = $value
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    $r.Font.Bold = $bold
    $r.Font.Size = $size
}

# ----------------------------------------------------------------------
# Values (unchanged cells are rewritten too, harmless, keeps this linear)
# ----------------------------------------------------------------------

# Row 1 - section titles
Format-Cell "A1" "hotel:" $false 14
Format-Cell "B1" "" $false 14
Format-Cell "D1" "room:" $false 14
Format-Cell "E1" "" $false 14

# Row 2 - column headers
Format-Cell "A2" "Cell name" $true
Format-Cell "B2" "Information" $true
Format-Cell "D2" "Cell name" $true
Format-Cell "E2" "Information" $true

# Row 3
Format-Cell "A3" "id_hotel"
Format-Cell "B3" "int, Not Null, Primary Key"
Format-Cell "D3" "id_room"
Format-Cell "E3" "int, Not Null, Primary Key"

# Row 4
Format-Cell "A4" "hotel_code"
Format-Cell "B4" "nvarchar(MAX), Not Null"
Format-Cell "D4" "tl_api_code"
Format-Cell "E4" "nvarchar(MAX), Not Null"

# Row 5
Format-Cell "A5" "login"
Format-Cell "B5" "nvarchar(MAX), Not Null"
Format-Cell "D5" "url"
Format-Cell "E5" "nvarchar(MAX), Not Null"

# Row 6 - "name" is the new field, inserted into the room table here
Format-Cell "A6" "password"
Format-Cell "B6" "nvarchar(MAX), Not Null"
Format-Cell "D6" "name"
Format-Cell "E6" "nvarchar(MAX), Not Null"

# Row 7 (new row) - "name" added at the bottom of the hotel table, and
# the room table's old row 6 (id_hotel / int, Not Null) moves here
Format-Cell "A7" "name"
Format-Cell "B7" "nvarchar(MAX), Not Null"
Format-Cell "D7" "id_hotel"
Format-Cell "E7" "int, Not Null"

$ws.Rows.Item(7).RowHeight = 29.4

# ----------------------------------------------------------------------
# Borders - every cell in the two tables gets medium left/right edges,
# thin edges between data rows, and a medium edge at the outer
# top/bottom boundary of each section.
# ----------------------------------------------------------------------

# Row 1 (title band) - medium top, no bottom (header row follows directly)
Set-CellBorders "A1" $xlMedium $xlThin   $xlMedium $null
Set-CellBorders "B1" $xlThin   $xlMedium $xlMedium $null
Set-CellBorders "D1" $xlMedium $xlThin   $xlMedium $null
Set-CellBorders "E1" $xlThin   $xlMedium $xlMedium $null

# Row 2 (column header band)
Set-CellBorders "A2" $xlMedium $xlMedium $xlMedium $null
Set-CellBorders "B2" $xlMedium $xlMedium $xlMedium $null
Set-CellBorders "D2" $xlMedium $xlMedium $xlMedium $xlMedium
Set-CellBorders "E2" $xlMedium $xlMedium $xlMedium $xlMedium

# Row 3
Set-CellBorders "A3" $xlMedium $xlMedium $xlMedium $xlThin
Set-CellBorders "B3" $xlMedium $xlMedium $xlMedium $xlThin
Set-CellBorders "D3" $xlMedium $xlMedium $null     $xlThin
Set-CellBorders "E3" $xlMedium $xlMedium $null     $xlThin

# Row 4
Set-CellBorders "A4" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "B4" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "D4" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "E4" $xlMedium $xlMedium $xlThin $xlThin

# Row 5
Set-CellBorders "A5" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "B5" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "D5" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "E5" $xlMedium $xlMedium $xlThin $xlThin

# Row 6
Set-CellBorders "A6" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "B6" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "D6" $xlMedium $xlMedium $xlThin $xlThin
Set-CellBorders "E6" $xlMedium $xlMedium $xlThin $xlThin

# Row 7 (new last row - medium bottom closes the table)
Set-CellBorders "A7" $xlMedium $xlMedium $xlThin $xlMedium
Set-CellBorders "B7" $xlMedium $xlMedium $xlThin $xlMedium
Set-CellBorders "D7" $xlMedium $xlMedium $null   $xlMedium
Set-CellBorders "E7" $xlMedium $xlMedium $null   $xlMedium

# ----------------------------------------------------------------------
# Sheet view bookkeeping to mirror the final state
# ----------------------------------------------------------------------
$ws.Range("B11").Select()
